# Update Papa (potato) price records for La Palmera de La Serena market.
# The data block in rows 387-472 is shifted down by two rows, two brand
# new records are inserted at rows 387-388, and the two previously-last
# records move into two brand new rows (473-474).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep a copy of the date-column number format (style index 2) before it
# gets applied to the two brand-new rows appended at the bottom.
$dateFmt = $ws.Cells.Item(386, 4).NumberFormat

# Row 387
$ws.Cells.Item(387, 4).Value = 44798
$ws.Cells.Item(387, 8).Value = 'Asterix'
$ws.Cells.Item(387, 9).Value = '1a (cosecha)'
$ws.Cells.Item(387, 10).Value = 2000
$ws.Cells.Item(387, 11).Value = 10000
$ws.Cells.Item(387, 12).Value = 10500
$ws.Cells.Item(387, 13).Value = 10250
$ws.Cells.Item(387, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(387, 16).Value = 410

# Row 388
$ws.Cells.Item(388, 4).Value = 44798
$ws.Cells.Item(388, 8).Value = 'Cardinal'
$ws.Cells.Item(388, 9).Value = '1a (cosecha)'
$ws.Cells.Item(388, 10).Value = 2000
$ws.Cells.Item(388, 11).Value = 12000
$ws.Cells.Item(388, 12).Value = 13000
$ws.Cells.Item(388, 13).Value = 12500
$ws.Cells.Item(388, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(388, 16).Value = 500

# Row 389
$ws.Cells.Item(389, 4).Value = 44342
$ws.Cells.Item(389, 8).Value = 'Asterix'
$ws.Cells.Item(389, 9).Value = '1a (cosecha)'
$ws.Cells.Item(389, 10).Value = 3000
$ws.Cells.Item(389, 11).Value = 7500
$ws.Cells.Item(389, 12).Value = 8000
$ws.Cells.Item(389, 13).Value = 7750
$ws.Cells.Item(389, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(389, 16).Value = 310

# Row 390
$ws.Cells.Item(390, 4).Value = 44342
$ws.Cells.Item(390, 8).Value = 'Cardinal'
$ws.Cells.Item(390, 9).Value = '1a nueva(o)'
$ws.Cells.Item(390, 10).Value = 2400
$ws.Cells.Item(390, 11).Value = 8000
$ws.Cells.Item(390, 12).Value = 9000
$ws.Cells.Item(390, 13).Value = 8500
$ws.Cells.Item(390, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(390, 16).Value = 340

# Row 391
$ws.Cells.Item(391, 4).Value = 44551
$ws.Cells.Item(391, 8).Value = 'Cardinal'
$ws.Cells.Item(391, 9).Value = '1a (cosecha)'
$ws.Cells.Item(391, 10).Value = 2400
$ws.Cells.Item(391, 11).Value = 12000
$ws.Cells.Item(391, 12).Value = 12500
$ws.Cells.Item(391, 13).Value = 12250
$ws.Cells.Item(391, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(391, 16).Value = 490

# Row 392
$ws.Cells.Item(392, 4).Value = 44648
$ws.Cells.Item(392, 8).Value = 'Asterix'
$ws.Cells.Item(392, 9).Value = '1a (cosecha)'
$ws.Cells.Item(392, 10).Value = 2460
$ws.Cells.Item(392, 11).Value = 8500
$ws.Cells.Item(392, 12).Value = 9000
$ws.Cells.Item(392, 13).Value = 8750
$ws.Cells.Item(392, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(392, 16).Value = 350

# Row 393
$ws.Cells.Item(393, 4).Value = 44291
$ws.Cells.Item(393, 8).Value = 'Cardinal'
$ws.Cells.Item(393, 9).Value = '1a (cosecha)'
$ws.Cells.Item(393, 10).Value = 2600
$ws.Cells.Item(393, 11).Value = 7500
$ws.Cells.Item(393, 12).Value = 8000
$ws.Cells.Item(393, 13).Value = 7750
$ws.Cells.Item(393, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(393, 16).Value = 310

# Row 394
$ws.Cells.Item(394, 4).Value = 44449
$ws.Cells.Item(394, 8).Value = 'Cardinal'
$ws.Cells.Item(394, 9).Value = '1a (cosecha)'
$ws.Cells.Item(394, 10).Value = 3000
$ws.Cells.Item(394, 11).Value = 12000
$ws.Cells.Item(394, 12).Value = 13000
$ws.Cells.Item(394, 13).Value = 12500
$ws.Cells.Item(394, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(394, 16).Value = 500

# Row 395
$ws.Cells.Item(395, 4).Value = 44585
$ws.Cells.Item(395, 8).Value = 'Asterix'
$ws.Cells.Item(395, 9).Value = '1a (cosecha)'
$ws.Cells.Item(395, 10).Value = 2500
$ws.Cells.Item(395, 11).Value = 9500
$ws.Cells.Item(395, 12).Value = 10000
$ws.Cells.Item(395, 13).Value = 9750
$ws.Cells.Item(395, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(395, 16).Value = 390

# Row 396
$ws.Cells.Item(396, 4).Value = 44376
$ws.Cells.Item(396, 8).Value = 'Cardinal'
$ws.Cells.Item(396, 9).Value = '1a nueva(o)'
$ws.Cells.Item(396, 10).Value = 2440
$ws.Cells.Item(396, 11).Value = 9000
$ws.Cells.Item(396, 12).Value = 9500
$ws.Cells.Item(396, 13).Value = 9250
$ws.Cells.Item(396, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(396, 16).Value = 370

# Row 397
$ws.Cells.Item(397, 4).Value = 44242
$ws.Cells.Item(397, 8).Value = 'Cardinal'
$ws.Cells.Item(397, 9).Value = '1a nueva(o)'
$ws.Cells.Item(397, 10).Value = 2500
$ws.Cells.Item(397, 11).Value = 9500
$ws.Cells.Item(397, 12).Value = 10000
$ws.Cells.Item(397, 13).Value = 9750
$ws.Cells.Item(397, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(397, 16).Value = 390

# Row 398
$ws.Cells.Item(398, 4).Value = 44391
$ws.Cells.Item(398, 8).Value = 'Cardinal'
$ws.Cells.Item(398, 9).Value = '1a nueva(o)'
$ws.Cells.Item(398, 10).Value = 3000
$ws.Cells.Item(398, 11).Value = 9500
$ws.Cells.Item(398, 12).Value = 10000
$ws.Cells.Item(398, 13).Value = 9750
$ws.Cells.Item(398, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(398, 16).Value = 390

# Row 399
$ws.Cells.Item(399, 4).Value = 44600
$ws.Cells.Item(399, 8).Value = 'Asterix'
$ws.Cells.Item(399, 9).Value = '1a (cosecha)'
$ws.Cells.Item(399, 10).Value = 2500
$ws.Cells.Item(399, 11).Value = 9500
$ws.Cells.Item(399, 12).Value = 10000
$ws.Cells.Item(399, 13).Value = 9750
$ws.Cells.Item(399, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(399, 16).Value = 390

# Row 400
$ws.Cells.Item(400, 4).Value = 44763
$ws.Cells.Item(400, 8).Value = 'Asterix'
$ws.Cells.Item(400, 9).Value = '1a (guarda)'
$ws.Cells.Item(400, 10).Value = 2000
$ws.Cells.Item(400, 11).Value = 11000
$ws.Cells.Item(400, 12).Value = 12000
$ws.Cells.Item(400, 13).Value = 11500
$ws.Cells.Item(400, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(400, 16).Value = 460

# Row 401
$ws.Cells.Item(401, 4).Value = 44371
$ws.Cells.Item(401, 8).Value = 'Asterix'
$ws.Cells.Item(401, 9).Value = '1a nueva(o)'
$ws.Cells.Item(401, 10).Value = 2400
$ws.Cells.Item(401, 11).Value = 8000
$ws.Cells.Item(401, 12).Value = 8500
$ws.Cells.Item(401, 13).Value = 8250
$ws.Cells.Item(401, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(401, 16).Value = 330

# Row 402
$ws.Cells.Item(402, 4).Value = 44371
$ws.Cells.Item(402, 8).Value = 'Cardinal'
$ws.Cells.Item(402, 9).Value = '1a nueva(o)'
$ws.Cells.Item(402, 10).Value = 2000
$ws.Cells.Item(402, 11).Value = 9000
$ws.Cells.Item(402, 12).Value = 9500
$ws.Cells.Item(402, 13).Value = 9250
$ws.Cells.Item(402, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(402, 16).Value = 370

# Row 403
$ws.Cells.Item(403, 4).Value = 44355
$ws.Cells.Item(403, 8).Value = 'Cardinal'
$ws.Cells.Item(403, 9).Value = '1a nueva(o)'
$ws.Cells.Item(403, 10).Value = 2400
$ws.Cells.Item(403, 11).Value = 8500
$ws.Cells.Item(403, 12).Value = 9000
$ws.Cells.Item(403, 13).Value = 8750
$ws.Cells.Item(403, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(403, 16).Value = 350

# Row 404
$ws.Cells.Item(404, 4).Value = 44579
$ws.Cells.Item(404, 8).Value = 'Asterix'
$ws.Cells.Item(404, 9).Value = '1a (cosecha)'
$ws.Cells.Item(404, 10).Value = 2400
$ws.Cells.Item(404, 11).Value = 10000
$ws.Cells.Item(404, 12).Value = 11000
$ws.Cells.Item(404, 13).Value = 10500
$ws.Cells.Item(404, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(404, 16).Value = 420

# Row 405
$ws.Cells.Item(405, 4).Value = 44558
$ws.Cells.Item(405, 8).Value = 'Asterix'
$ws.Cells.Item(405, 9).Value = '1a (cosecha)'
$ws.Cells.Item(405, 10).Value = 2200
$ws.Cells.Item(405, 11).Value = 10000
$ws.Cells.Item(405, 12).Value = 11000
$ws.Cells.Item(405, 13).Value = 10500
$ws.Cells.Item(405, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(405, 16).Value = 420

# Row 406
$ws.Cells.Item(406, 4).Value = 44558
$ws.Cells.Item(406, 8).Value = 'Cardinal'
$ws.Cells.Item(406, 9).Value = '1a (cosecha)'
$ws.Cells.Item(406, 10).Value = 2000
$ws.Cells.Item(406, 11).Value = 12000
$ws.Cells.Item(406, 12).Value = 13000
$ws.Cells.Item(406, 13).Value = 12500
$ws.Cells.Item(406, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(406, 16).Value = 500

# Row 407
$ws.Cells.Item(407, 4).Value = 44434
$ws.Cells.Item(407, 8).Value = 'Cardinal'
$ws.Cells.Item(407, 9).Value = '1a (cosecha)'
$ws.Cells.Item(407, 10).Value = 2000
$ws.Cells.Item(407, 11).Value = 9500
$ws.Cells.Item(407, 12).Value = 10000
$ws.Cells.Item(407, 13).Value = 9750
$ws.Cells.Item(407, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(407, 16).Value = 390

# Row 408
$ws.Cells.Item(408, 4).Value = 44462
$ws.Cells.Item(408, 8).Value = 'Cardinal'
$ws.Cells.Item(408, 9).Value = '1a (cosecha)'
$ws.Cells.Item(408, 10).Value = 2000
$ws.Cells.Item(408, 11).Value = 12000
$ws.Cells.Item(408, 12).Value = 13000
$ws.Cells.Item(408, 13).Value = 12500
$ws.Cells.Item(408, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(408, 16).Value = 500

# Row 409
$ws.Cells.Item(409, 4).Value = 44278
$ws.Cells.Item(409, 8).Value = 'Asterix'
$ws.Cells.Item(409, 9).Value = '1a (cosecha)'
$ws.Cells.Item(409, 10).Value = 2000
$ws.Cells.Item(409, 11).Value = 7500
$ws.Cells.Item(409, 12).Value = 8000
$ws.Cells.Item(409, 13).Value = 7750
$ws.Cells.Item(409, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(409, 16).Value = 310

# Row 410
$ws.Cells.Item(410, 4).Value = 44278
$ws.Cells.Item(410, 8).Value = 'Cardinal'
$ws.Cells.Item(410, 9).Value = '1a (cosecha)'
$ws.Cells.Item(410, 10).Value = 2400
$ws.Cells.Item(410, 11).Value = 8500
$ws.Cells.Item(410, 12).Value = 9000
$ws.Cells.Item(410, 13).Value = 8750
$ws.Cells.Item(410, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(410, 16).Value = 350

# Row 411
$ws.Cells.Item(411, 4).Value = 44442
$ws.Cells.Item(411, 8).Value = 'Cardinal'
$ws.Cells.Item(411, 9).Value = '1a (cosecha)'
$ws.Cells.Item(411, 10).Value = 3000
$ws.Cells.Item(411, 11).Value = 11500
$ws.Cells.Item(411, 12).Value = 12000
$ws.Cells.Item(411, 13).Value = 11750
$ws.Cells.Item(411, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(411, 16).Value = 470

# Row 412
$ws.Cells.Item(412, 4).Value = 44272
$ws.Cells.Item(412, 8).Value = 'Asterix'
$ws.Cells.Item(412, 9).Value = '1a (cosecha)'
$ws.Cells.Item(412, 10).Value = 2600
$ws.Cells.Item(412, 11).Value = 7500
$ws.Cells.Item(412, 12).Value = 8000
$ws.Cells.Item(412, 13).Value = 7750
$ws.Cells.Item(412, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(412, 16).Value = 310

# Row 413
$ws.Cells.Item(413, 4).Value = 44272
$ws.Cells.Item(413, 8).Value = 'Cardinal'
$ws.Cells.Item(413, 9).Value = '1a nueva(o)'
$ws.Cells.Item(413, 10).Value = 2400
$ws.Cells.Item(413, 11).Value = 8500
$ws.Cells.Item(413, 12).Value = 9000
$ws.Cells.Item(413, 13).Value = 8750
$ws.Cells.Item(413, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(413, 16).Value = 350

# Row 414
$ws.Cells.Item(414, 4).Value = 44238
$ws.Cells.Item(414, 8).Value = 'Cardinal'
$ws.Cells.Item(414, 9).Value = '1a nueva(o)'
$ws.Cells.Item(414, 10).Value = 2000
$ws.Cells.Item(414, 11).Value = 9500
$ws.Cells.Item(414, 12).Value = 10000
$ws.Cells.Item(414, 13).Value = 9750
$ws.Cells.Item(414, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(414, 16).Value = 390

# Row 415
$ws.Cells.Item(415, 4).Value = 44781
$ws.Cells.Item(415, 8).Value = 'Asterix'
$ws.Cells.Item(415, 9).Value = '1a (cosecha)'
$ws.Cells.Item(415, 10).Value = 2000
$ws.Cells.Item(415, 11).Value = 11000
$ws.Cells.Item(415, 12).Value = 12000
$ws.Cells.Item(415, 13).Value = 11500
$ws.Cells.Item(415, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(415, 16).Value = 460

# Row 416
$ws.Cells.Item(416, 4).Value = 44336
$ws.Cells.Item(416, 8).Value = 'Asterix'
$ws.Cells.Item(416, 9).Value = '1a (cosecha)'
$ws.Cells.Item(416, 10).Value = 2500
$ws.Cells.Item(416, 11).Value = 7000
$ws.Cells.Item(416, 12).Value = 7500
$ws.Cells.Item(416, 13).Value = 7250
$ws.Cells.Item(416, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(416, 16).Value = 290

# Row 417
$ws.Cells.Item(417, 4).Value = 44336
$ws.Cells.Item(417, 8).Value = 'Cardinal'
$ws.Cells.Item(417, 9).Value = '1a nueva(o)'
$ws.Cells.Item(417, 10).Value = 2000
$ws.Cells.Item(417, 11).Value = 8000
$ws.Cells.Item(417, 12).Value = 9000
$ws.Cells.Item(417, 13).Value = 8500
$ws.Cells.Item(417, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(417, 16).Value = 340

# Row 418
$ws.Cells.Item(418, 4).Value = 44343
$ws.Cells.Item(418, 8).Value = 'Asterix'
$ws.Cells.Item(418, 9).Value = '1a (cosecha)'
$ws.Cells.Item(418, 10).Value = 2500
$ws.Cells.Item(418, 11).Value = 7500
$ws.Cells.Item(418, 12).Value = 8000
$ws.Cells.Item(418, 13).Value = 7750
$ws.Cells.Item(418, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(418, 16).Value = 310

# Row 419
$ws.Cells.Item(419, 4).Value = 44343
$ws.Cells.Item(419, 8).Value = 'Cardinal'
$ws.Cells.Item(419, 9).Value = '1a nueva(o)'
$ws.Cells.Item(419, 10).Value = 2000
$ws.Cells.Item(419, 11).Value = 8000
$ws.Cells.Item(419, 12).Value = 9000
$ws.Cells.Item(419, 13).Value = 8500
$ws.Cells.Item(419, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(419, 16).Value = 340

# Row 420
$ws.Cells.Item(420, 4).Value = 44533
$ws.Cells.Item(420, 8).Value = 'Asterix'
$ws.Cells.Item(420, 9).Value = '1a nueva(o)'
$ws.Cells.Item(420, 10).Value = 2800
$ws.Cells.Item(420, 11).Value = 11500
$ws.Cells.Item(420, 12).Value = 12000
$ws.Cells.Item(420, 13).Value = 11750
$ws.Cells.Item(420, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(420, 16).Value = 470

# Row 421
$ws.Cells.Item(421, 4).Value = 44533
$ws.Cells.Item(421, 8).Value = 'Cardinal'
$ws.Cells.Item(421, 9).Value = '1a nueva(o)'
$ws.Cells.Item(421, 10).Value = 2800
$ws.Cells.Item(421, 11).Value = 12000
$ws.Cells.Item(421, 12).Value = 13000
$ws.Cells.Item(421, 13).Value = 12500
$ws.Cells.Item(421, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(421, 16).Value = 500

# Row 422
$ws.Cells.Item(422, 4).Value = 44365
$ws.Cells.Item(422, 8).Value = 'Asterix'
$ws.Cells.Item(422, 9).Value = '1a nueva(o)'
$ws.Cells.Item(422, 10).Value = 2900
$ws.Cells.Item(422, 11).Value = 8000
$ws.Cells.Item(422, 12).Value = 8500
$ws.Cells.Item(422, 13).Value = 8250
$ws.Cells.Item(422, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(422, 16).Value = 330

# Row 423
$ws.Cells.Item(423, 4).Value = 44454
$ws.Cells.Item(423, 8).Value = 'Cardinal'
$ws.Cells.Item(423, 9).Value = '1a (cosecha)'
$ws.Cells.Item(423, 10).Value = 2900
$ws.Cells.Item(423, 11).Value = 11500
$ws.Cells.Item(423, 12).Value = 12000
$ws.Cells.Item(423, 13).Value = 11750
$ws.Cells.Item(423, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(423, 16).Value = 470

# Row 424
$ws.Cells.Item(424, 4).Value = 44561
$ws.Cells.Item(424, 8).Value = 'Cardinal'
$ws.Cells.Item(424, 9).Value = '1a (cosecha)'
$ws.Cells.Item(424, 10).Value = 2900
$ws.Cells.Item(424, 11).Value = 12000
$ws.Cells.Item(424, 12).Value = 13000
$ws.Cells.Item(424, 13).Value = 12500
$ws.Cells.Item(424, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(424, 16).Value = 500

# Row 425
$ws.Cells.Item(425, 4).Value = 44421
$ws.Cells.Item(425, 8).Value = 'Cardinal'
$ws.Cells.Item(425, 9).Value = '1a (cosecha)'
$ws.Cells.Item(425, 10).Value = 3000
$ws.Cells.Item(425, 11).Value = 9500
$ws.Cells.Item(425, 12).Value = 10000
$ws.Cells.Item(425, 13).Value = 9750
$ws.Cells.Item(425, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(425, 16).Value = 390

# Row 426
$ws.Cells.Item(426, 4).Value = 44419
$ws.Cells.Item(426, 8).Value = 'Cardinal'
$ws.Cells.Item(426, 9).Value = '1a (cosecha)'
$ws.Cells.Item(426, 10).Value = 3100
$ws.Cells.Item(426, 11).Value = 9500
$ws.Cells.Item(426, 12).Value = 10000
$ws.Cells.Item(426, 13).Value = 9750
$ws.Cells.Item(426, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(426, 16).Value = 390

# Row 427
$ws.Cells.Item(427, 4).Value = 44669
$ws.Cells.Item(427, 8).Value = 'Asterix'
$ws.Cells.Item(427, 9).Value = '1a (cosecha)'
$ws.Cells.Item(427, 10).Value = 2400
$ws.Cells.Item(427, 11).Value = 8000
$ws.Cells.Item(427, 12).Value = 9000
$ws.Cells.Item(427, 13).Value = 8500
$ws.Cells.Item(427, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(427, 16).Value = 340

# Row 428
$ws.Cells.Item(428, 4).Value = 44699
$ws.Cells.Item(428, 8).Value = 'Asterix'
$ws.Cells.Item(428, 9).Value = '1a (cosecha)'
$ws.Cells.Item(428, 10).Value = 2000
$ws.Cells.Item(428, 11).Value = 9000
$ws.Cells.Item(428, 12).Value = 9500
$ws.Cells.Item(428, 13).Value = 9250
$ws.Cells.Item(428, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(428, 16).Value = 370

# Row 429
$ws.Cells.Item(429, 4).Value = 44636
$ws.Cells.Item(429, 8).Value = 'Asterix'
$ws.Cells.Item(429, 9).Value = '1a (cosecha)'
$ws.Cells.Item(429, 10).Value = 2000
$ws.Cells.Item(429, 11).Value = 8000
$ws.Cells.Item(429, 12).Value = 9000
$ws.Cells.Item(429, 13).Value = 8500
$ws.Cells.Item(429, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(429, 16).Value = 340

# Row 430
$ws.Cells.Item(430, 4).Value = 44405
$ws.Cells.Item(430, 8).Value = 'Cardinal'
$ws.Cells.Item(430, 9).Value = '1a nueva(o)'
$ws.Cells.Item(430, 10).Value = 3200
$ws.Cells.Item(430, 11).Value = 9500
$ws.Cells.Item(430, 12).Value = 10000
$ws.Cells.Item(430, 13).Value = 9750
$ws.Cells.Item(430, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(430, 16).Value = 390

# Row 431
$ws.Cells.Item(431, 4).Value = 44202
$ws.Cells.Item(431, 8).Value = 'Asterix'
$ws.Cells.Item(431, 9).Value = '1a nueva(o)'
$ws.Cells.Item(431, 10).Value = 2000
$ws.Cells.Item(431, 11).Value = 12500
$ws.Cells.Item(431, 12).Value = 13000
$ws.Cells.Item(431, 13).Value = 12750
$ws.Cells.Item(431, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(431, 16).Value = 510

# Row 432
$ws.Cells.Item(432, 4).Value = 44273
$ws.Cells.Item(432, 8).Value = 'Asterix'
$ws.Cells.Item(432, 9).Value = '1a (cosecha)'
$ws.Cells.Item(432, 10).Value = 2200
$ws.Cells.Item(432, 11).Value = 7000
$ws.Cells.Item(432, 12).Value = 7500
$ws.Cells.Item(432, 13).Value = 7250
$ws.Cells.Item(432, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(432, 16).Value = 290

# Row 433
$ws.Cells.Item(433, 4).Value = 44273
$ws.Cells.Item(433, 8).Value = 'Cardinal'
$ws.Cells.Item(433, 9).Value = '1a (cosecha)'
$ws.Cells.Item(433, 10).Value = 2000
$ws.Cells.Item(433, 11).Value = 8500
$ws.Cells.Item(433, 12).Value = 9000
$ws.Cells.Item(433, 13).Value = 8750
$ws.Cells.Item(433, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(433, 16).Value = 350

# Row 434
$ws.Cells.Item(434, 4).Value = 44777
$ws.Cells.Item(434, 8).Value = 'Asterix'
$ws.Cells.Item(434, 9).Value = '1a (cosecha)'
$ws.Cells.Item(434, 10).Value = 1900
$ws.Cells.Item(434, 11).Value = 11500
$ws.Cells.Item(434, 12).Value = 12000
$ws.Cells.Item(434, 13).Value = 11750
$ws.Cells.Item(434, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(434, 16).Value = 470

# Row 435
$ws.Cells.Item(435, 4).Value = 44159
$ws.Cells.Item(435, 8).Value = 'Cardinal'
$ws.Cells.Item(435, 9).Value = '1a (cosecha)'
$ws.Cells.Item(435, 10).Value = 2000
$ws.Cells.Item(435, 11).Value = 8500
$ws.Cells.Item(435, 12).Value = 9000
$ws.Cells.Item(435, 13).Value = 8750
$ws.Cells.Item(435, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(435, 16).Value = 350

# Row 436
$ws.Cells.Item(436, 4).Value = 44589
$ws.Cells.Item(436, 8).Value = 'Asterix'
$ws.Cells.Item(436, 9).Value = '1a (cosecha)'
$ws.Cells.Item(436, 10).Value = 2700
$ws.Cells.Item(436, 11).Value = 9500
$ws.Cells.Item(436, 12).Value = 10000
$ws.Cells.Item(436, 13).Value = 9750
$ws.Cells.Item(436, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(436, 16).Value = 390

# Row 437
$ws.Cells.Item(437, 4).Value = 44263
$ws.Cells.Item(437, 8).Value = 'Cardinal'
$ws.Cells.Item(437, 9).Value = '1a nueva(o)'
$ws.Cells.Item(437, 10).Value = 2400
$ws.Cells.Item(437, 11).Value = 8500
$ws.Cells.Item(437, 12).Value = 9000
$ws.Cells.Item(437, 13).Value = 8750
$ws.Cells.Item(437, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(437, 16).Value = 350

# Row 438
$ws.Cells.Item(438, 4).Value = 44309
$ws.Cells.Item(438, 8).Value = 'Asterix'
$ws.Cells.Item(438, 9).Value = '1a (cosecha)'
$ws.Cells.Item(438, 10).Value = 2800
$ws.Cells.Item(438, 11).Value = 7500
$ws.Cells.Item(438, 12).Value = 8000
$ws.Cells.Item(438, 13).Value = 7750
$ws.Cells.Item(438, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(438, 16).Value = 310

# Row 439
$ws.Cells.Item(439, 4).Value = 44771
$ws.Cells.Item(439, 8).Value = 'Asterix'
$ws.Cells.Item(439, 9).Value = '1a (guarda)'
$ws.Cells.Item(439, 10).Value = 2000
$ws.Cells.Item(439, 11).Value = 11000
$ws.Cells.Item(439, 12).Value = 12000
$ws.Cells.Item(439, 13).Value = 11500
$ws.Cells.Item(439, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(439, 16).Value = 460

# Row 440
$ws.Cells.Item(440, 4).Value = 44267
$ws.Cells.Item(440, 8).Value = 'Asterix'
$ws.Cells.Item(440, 9).Value = '1a (cosecha)'
$ws.Cells.Item(440, 10).Value = 2400
$ws.Cells.Item(440, 11).Value = 7500
$ws.Cells.Item(440, 12).Value = 8000
$ws.Cells.Item(440, 13).Value = 7750
$ws.Cells.Item(440, 15).Value = 'Región del Maule'
$ws.Cells.Item(440, 16).Value = 310

# Row 441
$ws.Cells.Item(441, 4).Value = 44267
$ws.Cells.Item(441, 8).Value = 'Cardinal'
$ws.Cells.Item(441, 9).Value = '1a nueva(o)'
$ws.Cells.Item(441, 10).Value = 2400
$ws.Cells.Item(441, 11).Value = 8000
$ws.Cells.Item(441, 12).Value = 9000
$ws.Cells.Item(441, 13).Value = 8500
$ws.Cells.Item(441, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(441, 16).Value = 340

# Row 442
$ws.Cells.Item(442, 4).Value = 44413
$ws.Cells.Item(442, 8).Value = 'Cardinal'
$ws.Cells.Item(442, 9).Value = '1a nueva(o)'
$ws.Cells.Item(442, 10).Value = 2000
$ws.Cells.Item(442, 11).Value = 9500
$ws.Cells.Item(442, 12).Value = 10000
$ws.Cells.Item(442, 13).Value = 9750
$ws.Cells.Item(442, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(442, 16).Value = 390

# Row 443
$ws.Cells.Item(443, 4).Value = 44328
$ws.Cells.Item(443, 8).Value = 'Asterix'
$ws.Cells.Item(443, 9).Value = '1a (cosecha)'
$ws.Cells.Item(443, 10).Value = 2900
$ws.Cells.Item(443, 11).Value = 7000
$ws.Cells.Item(443, 12).Value = 7500
$ws.Cells.Item(443, 13).Value = 7250
$ws.Cells.Item(443, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(443, 16).Value = 290

# Row 444
$ws.Cells.Item(444, 4).Value = 44515
$ws.Cells.Item(444, 8).Value = 'Cardinal'
$ws.Cells.Item(444, 9).Value = '1a nueva(o)'
$ws.Cells.Item(444, 10).Value = 2400
$ws.Cells.Item(444, 11).Value = 11500
$ws.Cells.Item(444, 12).Value = 12000
$ws.Cells.Item(444, 13).Value = 11750
$ws.Cells.Item(444, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(444, 16).Value = 470

# Row 445
$ws.Cells.Item(445, 4).Value = 44356
$ws.Cells.Item(445, 8).Value = 'Asterix'
$ws.Cells.Item(445, 9).Value = '1a (guarda)'
$ws.Cells.Item(445, 10).Value = 3000
$ws.Cells.Item(445, 11).Value = 8000
$ws.Cells.Item(445, 12).Value = 8500
$ws.Cells.Item(445, 13).Value = 8250
$ws.Cells.Item(445, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(445, 16).Value = 330

# Row 446
$ws.Cells.Item(446, 4).Value = 44356
$ws.Cells.Item(446, 8).Value = 'Cardinal'
$ws.Cells.Item(446, 9).Value = '1a nueva(o)'
$ws.Cells.Item(446, 10).Value = 2500
$ws.Cells.Item(446, 11).Value = 8500
$ws.Cells.Item(446, 12).Value = 9000
$ws.Cells.Item(446, 13).Value = 8750
$ws.Cells.Item(446, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(446, 16).Value = 350

# Row 447
$ws.Cells.Item(447, 4).Value = 44379
$ws.Cells.Item(447, 8).Value = 'Cardinal'
$ws.Cells.Item(447, 9).Value = '1a nueva(o)'
$ws.Cells.Item(447, 10).Value = 2800
$ws.Cells.Item(447, 11).Value = 9000
$ws.Cells.Item(447, 12).Value = 9500
$ws.Cells.Item(447, 13).Value = 9250
$ws.Cells.Item(447, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(447, 16).Value = 370

# Row 448
$ws.Cells.Item(448, 4).Value = 44322
$ws.Cells.Item(448, 8).Value = 'Asterix'
$ws.Cells.Item(448, 9).Value = '1a (cosecha)'
$ws.Cells.Item(448, 10).Value = 2520
$ws.Cells.Item(448, 11).Value = 7000
$ws.Cells.Item(448, 12).Value = 8000
$ws.Cells.Item(448, 13).Value = 7500
$ws.Cells.Item(448, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(448, 16).Value = 300

# Row 449
$ws.Cells.Item(449, 4).Value = 44221
$ws.Cells.Item(449, 8).Value = 'Asterix'
$ws.Cells.Item(449, 9).Value = '1a (cosecha)'
$ws.Cells.Item(449, 10).Value = 2600
$ws.Cells.Item(449, 11).Value = 9500
$ws.Cells.Item(449, 12).Value = 10000
$ws.Cells.Item(449, 13).Value = 9750
$ws.Cells.Item(449, 15).Value = 'Región del Maule'
$ws.Cells.Item(449, 16).Value = 390

# Row 450
$ws.Cells.Item(450, 4).Value = 44497
$ws.Cells.Item(450, 8).Value = 'Asterix'
$ws.Cells.Item(450, 9).Value = '1a nueva(o)'
$ws.Cells.Item(450, 10).Value = 2000
$ws.Cells.Item(450, 11).Value = 11500
$ws.Cells.Item(450, 12).Value = 12000
$ws.Cells.Item(450, 13).Value = 11750
$ws.Cells.Item(450, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(450, 16).Value = 470

# Row 451
$ws.Cells.Item(451, 4).Value = 44497
$ws.Cells.Item(451, 8).Value = 'Cardinal'
$ws.Cells.Item(451, 9).Value = '1a (cosecha)'
$ws.Cells.Item(451, 10).Value = 2400
$ws.Cells.Item(451, 11).Value = 12000
$ws.Cells.Item(451, 12).Value = 13000
$ws.Cells.Item(451, 13).Value = 12500
$ws.Cells.Item(451, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(451, 16).Value = 500

# Row 452
$ws.Cells.Item(452, 4).Value = 44782
$ws.Cells.Item(452, 8).Value = 'Asterix'
$ws.Cells.Item(452, 9).Value = '1a (cosecha)'
$ws.Cells.Item(452, 10).Value = 2000
$ws.Cells.Item(452, 11).Value = 11000
$ws.Cells.Item(452, 12).Value = 12000
$ws.Cells.Item(452, 13).Value = 11500
$ws.Cells.Item(452, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(452, 16).Value = 460

# Row 453
$ws.Cells.Item(453, 4).Value = 44435
$ws.Cells.Item(453, 8).Value = 'Cardinal'
$ws.Cells.Item(453, 9).Value = '1a (cosecha)'
$ws.Cells.Item(453, 10).Value = 11200
$ws.Cells.Item(453, 11).Value = 9000
$ws.Cells.Item(453, 12).Value = 10000
$ws.Cells.Item(453, 13).Value = 9612
$ws.Cells.Item(453, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(453, 16).Value = 384

# Row 454
$ws.Cells.Item(454, 4).Value = 44251
$ws.Cells.Item(454, 8).Value = 'Asterix'
$ws.Cells.Item(454, 9).Value = '1a nueva(o)'
$ws.Cells.Item(454, 10).Value = 2400
$ws.Cells.Item(454, 11).Value = 9000
$ws.Cells.Item(454, 12).Value = 10000
$ws.Cells.Item(454, 13).Value = 9500
$ws.Cells.Item(454, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(454, 16).Value = 380

# Row 455
$ws.Cells.Item(455, 4).Value = 44251
$ws.Cells.Item(455, 8).Value = 'Asterix'
$ws.Cells.Item(455, 9).Value = '1a nueva(o)'
$ws.Cells.Item(455, 10).Value = 3000
$ws.Cells.Item(455, 11).Value = 7500
$ws.Cells.Item(455, 12).Value = 8000
$ws.Cells.Item(455, 13).Value = 7750
$ws.Cells.Item(455, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(455, 16).Value = 310

# Row 456
$ws.Cells.Item(456, 4).Value = 44251
$ws.Cells.Item(456, 8).Value = 'Cardinal'
$ws.Cells.Item(456, 9).Value = '1a nueva(o)'
$ws.Cells.Item(456, 10).Value = 2600
$ws.Cells.Item(456, 11).Value = 9000
$ws.Cells.Item(456, 12).Value = 10000
$ws.Cells.Item(456, 13).Value = 9500
$ws.Cells.Item(456, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(456, 16).Value = 380

# Row 457
$ws.Cells.Item(457, 4).Value = 44319
$ws.Cells.Item(457, 8).Value = 'Asterix'
$ws.Cells.Item(457, 9).Value = '1a (cosecha)'
$ws.Cells.Item(457, 10).Value = 2700
$ws.Cells.Item(457, 11).Value = 7000
$ws.Cells.Item(457, 12).Value = 8000
$ws.Cells.Item(457, 13).Value = 7500
$ws.Cells.Item(457, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(457, 16).Value = 300

# Row 458
$ws.Cells.Item(458, 4).Value = 44344
$ws.Cells.Item(458, 8).Value = 'Asterix'
$ws.Cells.Item(458, 9).Value = '1a (cosecha)'
$ws.Cells.Item(458, 10).Value = 2860
$ws.Cells.Item(458, 11).Value = 7500
$ws.Cells.Item(458, 12).Value = 8000
$ws.Cells.Item(458, 13).Value = 7750
$ws.Cells.Item(458, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(458, 16).Value = 310

# Row 459
$ws.Cells.Item(459, 4).Value = 44232
$ws.Cells.Item(459, 8).Value = 'Asterix'
$ws.Cells.Item(459, 9).Value = '1a (cosecha)'
$ws.Cells.Item(459, 10).Value = 2000
$ws.Cells.Item(459, 11).Value = 8500
$ws.Cells.Item(459, 12).Value = 9000
$ws.Cells.Item(459, 13).Value = 8750
$ws.Cells.Item(459, 15).Value = 'Región del Maule'
$ws.Cells.Item(459, 16).Value = 350

# Row 460
$ws.Cells.Item(460, 4).Value = 44232
$ws.Cells.Item(460, 8).Value = 'Asterix'
$ws.Cells.Item(460, 9).Value = '1a nueva(o)'
$ws.Cells.Item(460, 10).Value = 2400
$ws.Cells.Item(460, 11).Value = 8500
$ws.Cells.Item(460, 12).Value = 9000
$ws.Cells.Item(460, 13).Value = 8750
$ws.Cells.Item(460, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(460, 16).Value = 350

# Row 461
$ws.Cells.Item(461, 4).Value = 44232
$ws.Cells.Item(461, 8).Value = 'Cardinal'
$ws.Cells.Item(461, 9).Value = '1a nueva(o)'
$ws.Cells.Item(461, 10).Value = 2000
$ws.Cells.Item(461, 11).Value = 9500
$ws.Cells.Item(461, 12).Value = 10000
$ws.Cells.Item(461, 13).Value = 9750
$ws.Cells.Item(461, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(461, 16).Value = 390

# Row 462
$ws.Cells.Item(462, 4).Value = 44455
$ws.Cells.Item(462, 8).Value = 'Cardinal'
$ws.Cells.Item(462, 9).Value = '1a (cosecha)'
$ws.Cells.Item(462, 10).Value = 2000
$ws.Cells.Item(462, 11).Value = 11500
$ws.Cells.Item(462, 12).Value = 12000
$ws.Cells.Item(462, 13).Value = 11750
$ws.Cells.Item(462, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(462, 16).Value = 470

# Row 463
$ws.Cells.Item(463, 4).Value = 44504
$ws.Cells.Item(463, 8).Value = 'Cardinal'
$ws.Cells.Item(463, 9).Value = '1a (cosecha)'
$ws.Cells.Item(463, 10).Value = 2400
$ws.Cells.Item(463, 11).Value = 12000
$ws.Cells.Item(463, 12).Value = 12500
$ws.Cells.Item(463, 13).Value = 12250
$ws.Cells.Item(463, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(463, 16).Value = 490

# Row 464
$ws.Cells.Item(464, 4).Value = 44484
$ws.Cells.Item(464, 8).Value = 'Cardinal'
$ws.Cells.Item(464, 9).Value = '1a (cosecha)'
$ws.Cells.Item(464, 10).Value = 2900
$ws.Cells.Item(464, 11).Value = 13000
$ws.Cells.Item(464, 12).Value = 14000
$ws.Cells.Item(464, 13).Value = 13500
$ws.Cells.Item(464, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(464, 16).Value = 540

# Row 465
$ws.Cells.Item(465, 4).Value = 44665
$ws.Cells.Item(465, 8).Value = 'Asterix'
$ws.Cells.Item(465, 9).Value = '1a (cosecha)'
$ws.Cells.Item(465, 10).Value = 2400
$ws.Cells.Item(465, 11).Value = 8000
$ws.Cells.Item(465, 12).Value = 9000
$ws.Cells.Item(465, 13).Value = 8500
$ws.Cells.Item(465, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(465, 16).Value = 340

# Row 466
$ws.Cells.Item(466, 4).Value = 44452
$ws.Cells.Item(466, 8).Value = 'Cardinal'
$ws.Cells.Item(466, 9).Value = '1a (cosecha)'
$ws.Cells.Item(466, 10).Value = 3000
$ws.Cells.Item(466, 11).Value = 11500
$ws.Cells.Item(466, 12).Value = 12000
$ws.Cells.Item(466, 13).Value = 11750
$ws.Cells.Item(466, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(466, 16).Value = 470

# Row 467
$ws.Cells.Item(467, 4).Value = 44510
$ws.Cells.Item(467, 8).Value = 'Asterix'
$ws.Cells.Item(467, 9).Value = '1a nueva(o)'
$ws.Cells.Item(467, 10).Value = 2500
$ws.Cells.Item(467, 11).Value = 11500
$ws.Cells.Item(467, 12).Value = 12000
$ws.Cells.Item(467, 13).Value = 11750
$ws.Cells.Item(467, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(467, 16).Value = 470

# Row 468
$ws.Cells.Item(468, 4).Value = 44510
$ws.Cells.Item(468, 8).Value = 'Cardinal'
$ws.Cells.Item(468, 9).Value = '1a (cosecha)'
$ws.Cells.Item(468, 10).Value = 2940
$ws.Cells.Item(468, 11).Value = 11500
$ws.Cells.Item(468, 12).Value = 12000
$ws.Cells.Item(468, 13).Value = 11750
$ws.Cells.Item(468, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(468, 16).Value = 470

# Row 469
$ws.Cells.Item(469, 4).Value = 44189
$ws.Cells.Item(469, 8).Value = 'Cardinal'
$ws.Cells.Item(469, 9).Value = '1a (cosecha)'
$ws.Cells.Item(469, 10).Value = 2000
$ws.Cells.Item(469, 11).Value = 12500
$ws.Cells.Item(469, 12).Value = 13000
$ws.Cells.Item(469, 13).Value = 12750
$ws.Cells.Item(469, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(469, 16).Value = 510

# Row 470
$ws.Cells.Item(470, 4).Value = 44701
$ws.Cells.Item(470, 8).Value = 'Asterix'
$ws.Cells.Item(470, 9).Value = '1a (cosecha)'
$ws.Cells.Item(470, 10).Value = 2520
$ws.Cells.Item(470, 11).Value = 8500
$ws.Cells.Item(470, 12).Value = 9000
$ws.Cells.Item(470, 13).Value = 8750
$ws.Cells.Item(470, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(470, 16).Value = 350

# Row 471
$ws.Cells.Item(471, 4).Value = 44516
$ws.Cells.Item(471, 8).Value = 'Cardinal'
$ws.Cells.Item(471, 9).Value = '1a nueva(o)'
$ws.Cells.Item(471, 10).Value = 2600
$ws.Cells.Item(471, 11).Value = 11500
$ws.Cells.Item(471, 12).Value = 12000
$ws.Cells.Item(471, 13).Value = 11750
$ws.Cells.Item(471, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(471, 16).Value = 470

# Row 472
$ws.Cells.Item(472, 4).Value = 44186
$ws.Cells.Item(472, 8).Value = 'Asterix'
$ws.Cells.Item(472, 9).Value = '1a nueva(o)'
$ws.Cells.Item(472, 10).Value = 2000
$ws.Cells.Item(472, 11).Value = 11500
$ws.Cells.Item(472, 12).Value = 12000
$ws.Cells.Item(472, 13).Value = 11750
$ws.Cells.Item(472, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(472, 16).Value = 470

# Row 473
$ws.Cells.Item(473, 4).Value = 44463
$ws.Cells.Item(473, 8).Value = 'Cardinal'
$ws.Cells.Item(473, 9).Value = '1a (cosecha)'
$ws.Cells.Item(473, 10).Value = 3000
$ws.Cells.Item(473, 11).Value = 12000
$ws.Cells.Item(473, 12).Value = 13000
$ws.Cells.Item(473, 13).Value = 12500
$ws.Cells.Item(473, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(473, 16).Value = 500

# Row 474
$ws.Cells.Item(474, 4).Value = 44382
$ws.Cells.Item(474, 8).Value = 'Cardinal'
$ws.Cells.Item(474, 9).Value = '1a (cosecha)'
$ws.Cells.Item(474, 10).Value = 3000
$ws.Cells.Item(474, 11).Value = 9000
$ws.Cells.Item(474, 12).Value = 9500
$ws.Cells.Item(474, 13).Value = 9250
$ws.Cells.Item(474, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(474, 16).Value = 370

# Rows 473-474 are entirely new rows; fill in the static columns too.
$ws.Cells.Item(473, 1).Value = 8
$ws.Cells.Item(473, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(473, 3).Value = 'Coquimbo'
$ws.Cells.Item(473, 5).Value = 4
$ws.Cells.Item(473, 6).Value = 100114001
$ws.Cells.Item(473, 7).Value = 'Papa'
$ws.Cells.Item(473, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(473, 17).Value = 25
$ws.Cells.Item(473, 18).Value = 'Hortaliza'
$ws.Cells.Item(473, 4).NumberFormat = $dateFmt

$ws.Cells.Item(474, 1).Value = 8
$ws.Cells.Item(474, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(474, 3).Value = 'Coquimbo'
$ws.Cells.Item(474, 5).Value = 4
$ws.Cells.Item(474, 6).Value = 100114001
$ws.Cells.Item(474, 7).Value = 'Papa'
$ws.Cells.Item(474, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(474, 17).Value = 25
$ws.Cells.Item(474, 18).Value = 'Hortaliza'
$ws.Cells.Item(474, 4).NumberFormat = $dateFmt

